$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Volume 30 Number 20" -> "Volume 30   Number  22" header text (A8)
$ws.Range("A8").Value = "Volume 30   Number  22"

# Update the "Report Covering the Week 5/15/2023 Through 5/21/2023" -> new dates (C9)
$ws.Range("C9").Value = "Report Covering the Week  5/29/2023  Through  6/4/2023"

# Update the weekly/28-day/YTD/2yr crime statistics table (rows 14-30)
$ws.Range("F14").Value = 8
$ws.Range("G14").Value = 11
$ws.Range("H14").Value = -27.272727272727
$ws.Range("I14").Value = 29
$ws.Range("J14").Value = 31
$ws.Range("K14").Value = -6.451612903225
$ws.Range("L14").Value = -17.142857142857
$ws.Range("M14").Value = -44.230769230769
$ws.Range("N14").Value = -84.974093264248
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 5
$ws.Range("E15").Value = -80
$ws.Range("F15").Value = 17
$ws.Range("G15").Value = 23
$ws.Range("H15").Value = -26.086956521739
$ws.Range("I15").Value = 94
$ws.Range("J15").Value = 103
$ws.Range("K15").Value = -8.737864077669
$ws.Range("L15").Value = 9.302325581395
$ws.Range("M15").Value = 1.075268817204
$ws.Range("N15").Value = -64.528301886792
$ws.Range("C16").Value = 50
$ws.Range("D16").Value = 38
$ws.Range("E16").Value = 31.578947368421
$ws.Range("F16").Value = 194
$ws.Range("G16").Value = 209
$ws.Range("H16").Value = -7.177033492822
$ws.Range("I16").Value = 999
$ws.Range("J16").Value = 1044
$ws.Range("K16").Value = -4.310344827586
$ws.Range("L16").Value = 23.791821561338
$ws.Range("M16").Value = -27.556200145032
$ws.Range("N16").Value = -85.254612546125
$ws.Range("C17").Value = 74
$ws.Range("D17").Value = 100
$ws.Range("E17").Value = -26
$ws.Range("F17").Value = 367
$ws.Range("G17").Value = 334
$ws.Range("H17").Value = 9.880239520958
$ws.Range("I17").Value = 1728
$ws.Range("J17").Value = 1641
$ws.Range("K17").Value = 5.301645338208
$ws.Range("L17").Value = 29.438202247191
$ws.Range("M17").Value = 30.809992429977
$ws.Range("N17").Value = -49.414519906323
$ws.Range("C18").Value = 31
$ws.Range("D18").Value = 43
$ws.Range("E18").Value = -27.906976744186
$ws.Range("F18").Value = 130
$ws.Range("G18").Value = 184
$ws.Range("H18").Value = -29.347826086956
$ws.Range("I18").Value = 850
$ws.Range("J18").Value = 1031
$ws.Range("K18").Value = -17.555771096023
$ws.Range("L18").Value = 11.695137976346
$ws.Range("M18").Value = -27.721088435374
$ws.Range("N18").Value = -82.412580177943
$ws.Range("C19").Value = 108
$ws.Range("D19").Value = 117
$ws.Range("E19").Value = -7.692307692307
$ws.Range("F19").Value = 441
$ws.Range("G19").Value = 456
$ws.Range("H19").Value = -3.289473684210
$ws.Range("I19").Value = 2384
$ws.Range("J19").Value = 2334
$ws.Range("K19").Value = 2.142245072836
$ws.Range("L19").Value = 37.327188940092
$ws.Range("M19").Value = 50.505050505050
$ws.Range("N19").Value = -8.903324417271
$ws.Range("C20").Value = 21
$ws.Range("D20").Value = 25
$ws.Range("E20").Value = -16
$ws.Range("F20").Value = 151
$ws.Range("G20").Value = 124
$ws.Range("H20").Value = 21.774193548387
$ws.Range("I20").Value = 713
$ws.Range("J20").Value = 717
$ws.Range("K20").Value = -0.557880055788
$ws.Range("L20").Value = 32.774674115456
$ws.Range("M20").Value = 24
$ws.Range("N20").Value = -82.537349987754
$ws.Range("C21").Value = 286
$ws.Range("D21").Value = 331
$ws.Range("E21").Value = -13.595166163142
$ws.Range("F21").Value = 1308
$ws.Range("G21").Value = 1341
$ws.Range("H21").Value = -2.460850111856
$ws.Range("I21").Value = 6797
$ws.Range("J21").Value = 6901
$ws.Range("K21").Value = -1.507027966961
$ws.Range("L21").Value = 28.317915801397
$ws.Range("M21").Value = 9.983818770226
$ws.Range("N21").Value = -69.358038048868
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 8
$ws.Range("E22").Value = -62.5
$ws.Range("F22").Value = 19
$ws.Range("G22").Value = 29
$ws.Range("H22").Value = -34.482758620689
$ws.Range("I22").Value = 125
$ws.Range("J22").Value = 160
$ws.Range("K22").Value = -21.875
$ws.Range("L22").Value = 15.740740740740
$ws.Range("M22").Value = -33.155080213903
$ws.Range("C23").Value = 36
$ws.Range("D23").Value = 26
$ws.Range("E23").Value = 38.461538461538
$ws.Range("F23").Value = 120
$ws.Range("G23").Value = 121
$ws.Range("H23").Value = -0.826446280991
$ws.Range("I23").Value = 667
$ws.Range("J23").Value = 615
$ws.Range("K23").Value = 8.455284552845
$ws.Range("L23").Value = 15.198618307426
$ws.Range("M23").Value = 50.904977375565
$ws.Range("C24").Value = 242
$ws.Range("D24").Value = 234
$ws.Range("E24").Value = 3.418803418803
$ws.Range("F24").Value = 960
$ws.Range("G24").Value = 1094
$ws.Range("H24").Value = -12.248628884826
$ws.Range("I24").Value = 5145
$ws.Range("J24").Value = 5363
$ws.Range("K24").Value = -4.064889054633
$ws.Range("L24").Value = 26.412776412776
$ws.Range("M24").Value = 28.304239401496
$ws.Range("C25").Value = 136
$ws.Range("D25").Value = 175
$ws.Range("E25").Value = -22.285714285714
$ws.Range("F25").Value = 517
$ws.Range("G25").Value = 554
$ws.Range("H25").Value = -6.678700361010
$ws.Range("I25").Value = 2519
$ws.Range("J25").Value = 2570
$ws.Range("K25").Value = -1.984435797665
$ws.Range("L25").Value = 41.995490417136
$ws.Range("M25").Value = -22.966360856269
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = -42.857142857142
$ws.Range("F26").Value = 31
$ws.Range("G26").Value = 35
$ws.Range("H26").Value = -11.428571428571
$ws.Range("I26").Value = 141
$ws.Range("J26").Value = 158
$ws.Range("K26").Value = -10.759493670886
$ws.Range("L26").Value = -8.441558441558
$ws.Range("C27").Value = 17
$ws.Range("D27").Value = 16
$ws.Range("E27").Value = 6.25
$ws.Range("G27").Value = 55
$ws.Range("H27").Value = -1.818181818181
$ws.Range("I27").Value = 257
$ws.Range("J27").Value = 252
$ws.Range("K27").Value = 1.984126984126
$ws.Range("L27").Value = -5.860805860805
$ws.Range("D28").Value = 15
$ws.Range("E28").Value = -66.666666666666
$ws.Range("G28").Value = 35
$ws.Range("H28").Value = -42.857142857142
$ws.Range("I28").Value = 98
$ws.Range("J28").Value = 122
$ws.Range("K28").Value = -19.672131147541
$ws.Range("L28").Value = -35.526315789473
$ws.Range("M28").Value = -45.856353591160
$ws.Range("N28").Value = -87.594936708860
$ws.Range("D29").Value = 12
$ws.Range("E29").Value = -58.333333333333
$ws.Range("F29").Value = 19
$ws.Range("G29").Value = 28
$ws.Range("H29").Value = -32.142857142857
$ws.Range("I29").Value = 85
$ws.Range("J29").Value = 105
$ws.Range("K29").Value = -19.047619047619
$ws.Range("L29").Value = -34.615384615384
$ws.Range("M29").Value = -40.559440559440
$ws.Range("N29").Value = -88.194444444444
$ws.Range("D30").Value = 4
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 9
$ws.Range("H30").Value = -77.777777777777
$ws.Range("J30").Value = 32
$ws.Range("K30").Value = -9.375
$ws.Range("L30").Value = 16
